$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Column A holds an empty-but-present string cell in this sheet (same as
# A2:A12 above it). A leading apostrophe forces Excel to commit a literal
# (non-blank) text cell instead of discarding an empty assignment; the
# style is then reset to "Normal" so no extra number-format style sticks
# to the cell (matches the plain, unstyled cells used throughout the rest
# of the sheet).
$cA = $ws.Cells.Item($row, 1)
$cA.Value = "'"
$cA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "احمد"

# Column C values look numeric ("222") but the sheet stores quantities as
# text everywhere else, so force text entry the same way a user typing
# '222 into the cell would, then strip the resulting style back to Normal.
$cC = $ws.Cells.Item($row, 3)
$cC.Value = "'222"
$cC.Style = "Normal"

$ws.Cells.Item($row, 4).Value = "الجزائري"
$ws.Cells.Item($row, 5).Value = "الرحلة 2"
$ws.Cells.Item($row, 6).Value = "C3"
$ws.Cells.Item($row, 7).Value = "WCK"
$ws.Cells.Item($row, 8).Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠٢:٤٠:٥٥ م"
